$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell text content ---

# Row 3 (C3): insert "加大worker数，" into item 4 of the optimization notes
$ws.Range("C3").Value = "1. 分析存储过程spMaterialByBatchID：`n  - 尝试在子查询加过滤条件，速度明显提升（30S->5S）；缺陷：在前端不输入过滤条件时，子查询无法过滤数据造成读取次数千万级以上；建议：在不输入查询条件时，仅返回有限的记录数`n  - 尝试在一些连接次数多的表建立索引并改用like xxx%，速度提升40%（5S->2S）；缺陷：只能匹配结尾；`n2. 分析视图vm_BatchStoreQty，尝试减少聚集计算的左外连接数（由4个减为2个）但效果不明显`n3. 压力测试存储过程并监测所占用的进程数和worker数，模拟触发“查询超时已过期”错误`n4. 分析数据连接数，查询等待时间及并行度设置，通过更改配置，加大worker数，调试结果基本无区别。`n5. 压力测试优化过的存储过程，模拟100个客户同时领料查询并根据worker数量控制客户端可以真正进入查询或等待，通过进程的干预暂未触发“查询超时已过期”错误"
$ws.Rows(3).RowHeight = 148.5

# Row 4 (C4): append a new item 3 to the PDA outbound analysis notes
$ws.Range("C4").Value = "1. 分析存储过程sp_pdaOutCaseCommit`n建议： `n   - SET SACT_ABORT ON`n   - C#应用程序加错误捕捉和处理语句（处理超时异常）`n2. 分析“General Network Error,Check your Network Documentation”`n3. 在测试环境将出库存储过程设置一定的等待时间，PDA提交数据后断开WIFI，可模拟出实际网络出问题后的错误。所有更改都已回滚，数据库并无死锁产生，网络恢复后用户可再次提交数据。`n建议：解决该问题仅能从网络线路入手"
$ws.Rows(4).RowHeight = 108

# Row 5 (B5): rename the fourth feature entry from 报表 to SQL优化
$ws.Range("B5").Value = "SQL优化"

# --- Add new rows 6, 7, 8 with the new SQL optimization entries ---

$ws.Range("B6").Value = "SQL11"
$ws.Range("C6").Value = "返回数据723024，用时4分钟`n"
$ws.Range("C6").WrapText = $true
$ws.Range("C6").VerticalAlignment = -4160
$ws.Rows(6).RowHeight = 27

$ws.Range("B7").Value = "SQL12"
$ws.Range("C7").Value = "调试前存储过程中的查询部分运行时间为1:30左右，记录数目13748， insert 13748条数据`n分析：瓶颈在用户自定义函数GetstockPasterColorBatchCrokid，已验证`n建议：改成表连接来获取批号和缸号（咨询过原作者夏工，对方声称该自动以函数没问题并且结果是一对多，即订单号+物料编号+颜色编号可能对多个批号）`n另发现GetstockPasterColorBatchCrokid一个bug，会造成批号或者缸号第一位丢失，已告诉夏工修复该bug`n由于该存储过程为作业调用，一天仅执行两次切在低谷期执行，对系统性能不会造成太大影响，可考虑不优化"
$ws.Range("C7").WrapText = $true
$ws.Range("C7").VerticalAlignment = -4160
$ws.Rows(7).RowHeight = 108

$ws.Range("B8").Value = "SQL13"
$ws.Range("C8").Value = "调试前运行时间为5分钟左右`n建议：`n- rev_detail增加非聚集不唯一索引idx_tmprecdate，运行时间提升为5秒"
$ws.Range("C8").WrapText = $true
$ws.Range("C8").VerticalAlignment = -4160
$ws.Rows(8).RowHeight = 40.5

# --- Update selection to match target view state ---
$ws.Range("C6").Select()
